$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Stocks": update Rendement moyen / Risque values, add Isin column (D)
# ---------------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("Stocks")

$stocks.Range("C1").Copy()
$stocks.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$stocks.Range("D1").Value = "Isin"

$stocksData = @(
    @{ Row = 2;  B = 0.0009372363294920203; C = 0.01806883302664484; D = "FR0000121014" },
    @{ Row = 3;  B = 0.0006288009129108965; C = 0.01907975781065794; D = "FR0000120271" },
    @{ Row = 4;  B = 0.0003693315763554386; C = 0.01405945875487118; D = "FR0000120578" },
    @{ Row = 5;  B = 0.0007411005344758494; C = 0.01477996623511113; D = "FR0000120321" },
    @{ Row = 6;  B = 0.0008204089533400808; C = 0.01751968796903483; D = "FR0000121972" },
    @{ Row = 7;  B = 0.0006224884847651957; C = 0.01297528881552609; D = "FR0000120073" },
    @{ Row = 8;  B = 0.0006954024960210487; C = 0.02515388240841759; D = "NL0000235190" },
    @{ Row = 9;  B = 0.0003668726422083456; C = 0.02114005613800019; D = "FR0000131104" },
    @{ Row = 10; B = 0.001149469027186647;  C = 0.01641119594322464; D = "FR0000052292" },
    @{ Row = 11; B = 0.0004851770110929554; C = 0.01615083025323356; D = "FR0000121667" }
)

foreach ($entry in $stocksData) {
    $r = $entry.Row
    $stocks.Cells.Item($r, 2).Value = $entry.B
    $stocks.Cells.Item($r, 3).Value = $entry.C
    $stocks.Cells.Item($r, 4).Value = $entry.D
}

# ---------------------------------------------------------------------------
# Sheet "Index": update Rendement moyen / Risque values
# ---------------------------------------------------------------------------
$index = $wb.Worksheets.Item("Index")

$indexData = @(
    @{ Row = 2; B = -0.0000128447583177529; C = 0.01201184060693191 },
    @{ Row = 3; B = 0.0003072664386585315;  C = 0.01279721437434423 },
    @{ Row = 4; B = 0.0002772623732336768;  C = 0.01237062834888296 },
    @{ Row = 5; B = 0.0002144727224379599;  C = 0.01270016000634717 },
    @{ Row = 6; B = 0.0001308948326562523;  C = 0.009590478528196261 }
)

foreach ($entry in $indexData) {
    $r = $entry.Row
    $index.Cells.Item($r, 2).Value = $entry.B
    $index.Cells.Item($r, 3).Value = $entry.C
}
